$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 46.96651666666666
$ws.Range("H2").Value = 140.89955
$ws.Range("I2").Value = 0.5808027674561179
$ws.Range("J2").Value = 0.5808027674561179
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.788187000000001
$ws.Range("N2").Value = 29.364561
$ws.Range("O2").Value = 0.1731793198378281
$ws.Range("P2").Value = 0.1731793198378281
$ws.Range("Q2").Value = 459.7170478719499
$ws.Range("R2").Value = 4137.453430847549
$ws.Range("S2").Value = 0.1005830282279788
$ws.Range("T2").Value = 0.1005830282279788
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 46.96651666666666
$ws.Range("H3").Value = 140.89955
$ws.Range("I3").Value = 0.5808027674561179
$ws.Range("J3").Value = 0.5808027674561179
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.05628333333334
$ws.Range("N3").Value = 54.16885000000001
$ws.Range("O3").Value = 0.3194641527042525
$ws.Range("P3").Value = 0.3194641527042525
$ws.Range("Q3").Value = 848.0407321130555
$ws.Range("R3").Value = 7632.3665890175
$ws.Range("S3").Value = 0.1855456639936537
$ws.Range("T3").Value = 0.1855456639936537
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 46.96651666666666
$ws.Range("H4").Value = 140.89955
$ws.Range("I4").Value = 0.5808027674561179
$ws.Range("J4").Value = 0.5808027674561179
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.510488
$ws.Range("N4").Value = 43.531464
$ws.Range("O4").Value = 0.2567295089841425
$ws.Range("P4").Value = 0.2567295089841425
$ws.Range("Q4").Value = 681.5070764934666
$ws.Range("R4").Value = 6133.563688441199
$ws.Range("S4").Value = 0.1491092093056403
$ws.Range("T4").Value = 0.1491092093056403
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 46.96651666666666
$ws.Range("H5").Value = 140.89955
$ws.Range("I5").Value = 0.5808027674561179
$ws.Range("J5").Value = 0.5808027674561179
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.165572
$ws.Range("N5").Value = 42.496716
$ws.Range("O5").Value = 0.2506270184737769
$ws.Range("P5").Value = 0.2506270184737769
$ws.Range("Q5").Value = 665.3075734308665
$ws.Range("R5").Value = 5987.768160877799
$ws.Range("S5").Value = 0.1455648659288452
$ws.Range("T5").Value = 0.1455648659288452
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.34807866666667
$ws.Range("H6").Value = 43.04423600000001
$ws.Range("I6").Value = 0.1774328689611448
$ws.Range("J6").Value = 0.1774328689611448
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.788187000000001
$ws.Range("N6").Value = 29.364561
$ws.Range("O6").Value = 0.1731793198378281
$ws.Range("P6").Value = 0.1731793198378281
$ws.Range("Q6").Value = 140.441677080044
$ws.Range("R6").Value = 1263.975093720396
$ws.Range("S6").Value = 0.03072770356356554
$ws.Range("T6").Value = 0.03072770356356554
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.34807866666667
$ws.Range("H7").Value = 43.04423600000001
$ws.Range("I7").Value = 0.1774328689611448
$ws.Range("J7").Value = 0.1774328689611448
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.05628333333334
$ws.Range("N7").Value = 54.16885000000001
$ws.Range("O7").Value = 0.3194641527042525
$ws.Range("P7").Value = 0.3194641527042525
$ws.Range("Q7").Value = 259.072973694289
$ws.Range("R7").Value = 2331.6567632486
$ws.Range("S7").Value = 0.05668344114455679
$ws.Range("T7").Value = 0.05668344114455678
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.34807866666667
$ws.Range("H8").Value = 43.04423600000001
$ws.Range("I8").Value = 0.1774328689611448
$ws.Range("J8").Value = 0.1774328689611448
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.510488
$ws.Range("N8").Value = 43.531464
$ws.Range("O8").Value = 0.2567295089841425
$ws.Range("P8").Value = 0.2567295089841425
$ws.Range("Q8").Value = 208.1976233157227
$ws.Range("R8").Value = 1873.778609841504
$ws.Range("S8").Value = 0.04555225332604241
$ws.Range("T8").Value = 0.0455522533260424
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.34807866666667
$ws.Range("H9").Value = 43.04423600000001
$ws.Range("I9").Value = 0.1774328689611448
$ws.Range("J9").Value = 0.1774328689611448
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.165572
$ws.Range("N9").Value = 42.496716
$ws.Range("O9").Value = 0.2506270184737769
$ws.Range("P9").Value = 0.2506270184737769
$ws.Range("Q9").Value = 203.2487414143307
$ws.Range("R9").Value = 1829.238672728976
$ws.Range("S9").Value = 0.04446947092698007
$ws.Range("T9").Value = 0.04446947092698007
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.746361333333333
$ws.Range("H10").Value = 5.239084
$ws.Range("I10").Value = 0.02159605538935411
$ws.Range("J10").Value = 0.02159605538935411
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.788187000000001
$ws.Range("N10").Value = 29.364561
$ws.Range("O10").Value = 0.1731793198378281
$ws.Range("P10").Value = 0.1731793198378281
$ws.Range("Q10").Value = 17.093711300236
$ws.Range("R10").Value = 153.843401702124
$ws.Range("S10").Value = 0.003739990183508408
$ws.Range("T10").Value = 0.003739990183508407
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.746361333333333
$ws.Range("H11").Value = 5.239084
$ws.Range("I11").Value = 0.02159605538935411
$ws.Range("J11").Value = 0.02159605538935411
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 18.05628333333334
$ws.Range("N11").Value = 54.16885000000001
$ws.Range("O11").Value = 0.3194641527042525
$ws.Range("P11").Value = 0.3194641527042525
$ws.Range("Q11").Value = 31.53279503704445
$ws.Range("R11").Value = 283.7951553334
$ws.Range("S11").Value = 0.006899165536714118
$ws.Range("T11").Value = 0.006899165536714115
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.746361333333333
$ws.Range("H12").Value = 5.239084
$ws.Range("I12").Value = 0.02159605538935411
$ws.Range("J12").Value = 0.02159605538935411
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 14.510488
$ws.Range("N12").Value = 43.531464
$ws.Range("O12").Value = 0.2567295089841425
$ws.Range("P12").Value = 0.2567295089841425
$ws.Range("Q12").Value = 25.34055517099734
$ws.Range("R12").Value = 228.064996538976
$ws.Range("S12").Value = 0.005544344696103226
$ws.Range("T12").Value = 0.005544344696103225
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.746361333333333
$ws.Range("H13").Value = 5.239084
$ws.Range("I13").Value = 0.02159605538935411
$ws.Range("J13").Value = 0.02159605538935411
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 14.165572
$ws.Range("N13").Value = 42.496716
$ws.Range("O13").Value = 0.2506270184737769
$ws.Range("P13").Value = 0.2506270184737769
$ws.Range("Q13").Value = 24.73820720534933
$ws.Range("R13").Value = 222.643864848144
$ws.Range("S13").Value = 0.005412554973028362
$ws.Range("T13").Value = 0.005412554973028361
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.80387266666667
$ws.Range("H14").Value = 53.411618
$ws.Range("I14").Value = 0.2201683081933832
$ws.Range("J14").Value = 0.2201683081933832
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 9.788187000000001
$ws.Range("N14").Value = 29.364561
$ws.Range("O14").Value = 0.1731793198378281
$ws.Range("P14").Value = 0.1731793198378281
$ws.Range("Q14").Value = 174.267634985522
$ws.Range("R14").Value = 1568.408714869698
$ws.Range("S14").Value = 0.03812859786277543
$ws.Range("T14").Value = 0.03812859786277543
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.80387266666667
$ws.Range("H15").Value = 53.411618
$ws.Range("I15").Value = 0.2201683081933832
$ws.Range("J15").Value = 0.2201683081933832
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 18.05628333333334
$ws.Range("N15").Value = 54.16885000000001
$ws.Range("O15").Value = 0.3194641527042525
$ws.Range("P15").Value = 0.3194641527042525
$ws.Range("Q15").Value = 321.4717692999223
$ws.Range("R15").Value = 2893.2459236993
$ws.Range("S15").Value = 0.07033588202932792
$ws.Range("T15").Value = 0.0703358820293279
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.80387266666667
$ws.Range("H16").Value = 53.411618
$ws.Range("I16").Value = 0.2201683081933832
$ws.Range("J16").Value = 0.2201683081933832
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 14.510488
$ws.Range("N16").Value = 43.531464
$ws.Range("O16").Value = 0.2567295089841425
$ws.Range("P16").Value = 0.2567295089841425
$ws.Range("Q16").Value = 258.3428806831947
$ws.Range("R16").Value = 2325.085926148752
$ws.Range("S16").Value = 0.05652370165635665
$ws.Range("T16").Value = 0.05652370165635664
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.80387266666667
$ws.Range("H17").Value = 53.411618
$ws.Range("I17").Value = 0.2201683081933832
$ws.Range("J17").Value = 0.2201683081933832
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 14.165572
$ws.Range("N17").Value = 42.496716
$ws.Range("O17").Value = 0.2506270184737769
$ws.Range("P17").Value = 0.2506270184737769
$ws.Range("Q17").Value = 252.2020401384987
$ws.Range("R17").Value = 2269.818361246488
$ws.Range("S17").Value = 0.05518012664492326
$ws.Range("T17").Value = 0.05518012664492326
